$d = $word.ActiveDocument

# 1. Update the title on the cover page.
$d.Content.Find.Execute("Project Report", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "Fruits Classification Using Convolutional Neural Network", 2) | Out-Null

# 2. Remove the last 4 blank paragraphs of the cover-page spacer run
#    (right before the "Introduction" heading), deleting from the
#    bottom up so earlier indices stay valid.
$introRng = $d.Content
$introRng.Find.Execute("Introduction") | Out-Null
$introParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($introRng.Start -ge $p.Range.Start -and $introRng.Start -lt $p.Range.End) {
        $introParaIndex = $i
        break
    }
}
for ($n = 1; $n -le 4; $n++) {
    $idx = $introParaIndex - $n
    $d.Paragraphs.Item($idx).Range.Delete() | Out-Null
}

# 3. Move the "_GoBack" bookmark from the blank paragraph after
#    "Experimental Setup" to right after the "Deep Learning Neural
#    Network" heading text. A bookmark collapsed exactly at a
#    paragraph's end position can't be created directly, so insert a
#    marker character, bookmark it, then delete the marker - this
#    leaves the bookmark collapsed in the right spot.
$rng = $d.Content
$rng.Find.Execute("Deep Learning Neural Network") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
$rng.Text = ""
